$wb = $excel.ActiveWorkbook

# Update the unit cost value on the economical_params sheet from 140 €/kg to 0.23 €/g
$wsParams = $wb.Worksheets.Item("economical_params")
$cell = $wsParams.Range("B2")
$cell.Value = 0.23

# Re-apply the cell's alignment so the style record collapses onto the
# existing (non number-format-flagged) style instead of keeping the
# redundant one around.
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108

# Move selection to B8 as recorded in the saved workbook view
$wsParams.Range("B8").Select()

$wb.Save()
